$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: Base Plate CAD file / cost of machining ---
$baseUrl = "https://github.com/UCSD-E4E/3d-visualization-system/blob/master/CAD/Base%20plate.SLDPRT"
$ws.Range("E15").Value = $baseUrl
$ws.Hyperlinks.Add($ws.Range("E15"), $baseUrl)
$ws.Range("E15").Style = "Hyperlink"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 175
$ws.Range("H15").Formula = "=F15*G15"

# --- Row 16: Enclosure tube holder CAD file / cost of machining ---
$holderUrl = "https://github.com/UCSD-E4E/3d-visualization-system/blob/master/CAD/enclosure%20tube%20holder.SLDPRT"
$ws.Range("E16").Value = $holderUrl
$ws.Hyperlinks.Add($ws.Range("E16"), $holderUrl)
$ws.Range("E16").Style = "Hyperlink"
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 125
$ws.Range("H16").Formula = "=F16*G16"

# --- New "Notes" column header ---
$ws.Range("I1").Value = "Notes"

# --- Notes for machine-shop screws that are included in machining cost ---
$ws.Range("I17").Value = "(included in machining cost)"
$ws.Range("I18").Value = "(included in machining cost)"

# --- Restore the active selection to where the editor left off ---
[void]$ws.Range("I19").Select()
